# Correct status name labels in column B (statut_label) and column C (statut_name)
# - "bleu" -> "noir"
# - "pas de résultat ni de publication" -> "pas de résultat postés ni publiés"
# - "résultat et / ou publication posté dans les 12 mois" -> "résultat postés ou publiés dans les 12 mois"
# - "résultat et / ou publication posté dans les 36 mois" -> "résultat postés ou publiés dans les 36 mois"
# - "résultat et / ou publication posté" -> "résultat postés ou publiés"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$replacements = @{
    "bleu" = "noir"
    "pas de résultat ni de publication" = "pas de résultat postés ni publiés"
    "résultat et / ou publication posté dans les 12 mois" = "résultat postés ou publiés dans les 12 mois"
    "résultat et / ou publication posté dans les 36 mois" = "résultat postés ou publiés dans les 36 mois"
    "résultat et / ou publication posté" = "résultat postés ou publiés"
}

$usedRange = $ws.UsedRange
$rowCount = $usedRange.Rows.Count
$colCount = $usedRange.Columns.Count

for ($r = 1; $r -le $rowCount; $r++) {
    for ($c = 1; $c -le $colCount; $c++) {
        $cell = $ws.Cells.Item($r, $c)
        $val = $cell.Value2
        if ($null -ne $val -and $replacements.ContainsKey($val)) {
            $cell.Value = $replacements[$val]
        }
    }
}
